# Apogee : mise a jour 4/10/18
# "Nationalite" sheet: insert a new country row, BELIZE (code 429), right
# after GUYANAIS(E) (code 428). Every following row shifts down by one,
# which pushes the last data row (999 / INCONNUE) from row 181 to row 182,
# and the trailing blank spacer row from row 182 to row 183.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nationalité")

$firstRow = 166
$lastRow  = 181

# --- Grow the table by one row --------------------------------------------
$spacerRow = $lastRow + 1
$newSpacerRow = $lastRow + 2

# New row 183 becomes the trailing blank spacer row: same height as the
# current spacer row (182), but left cell-free (no formatting copy, so no
# stray empty cells get created - matches the original untouched spacer
# row's shape).
$ws.Rows.Item($newSpacerRow).RowHeight = $ws.Rows.Item($spacerRow).RowHeight

# New row 182 becomes an ordinary data row, matching the alternating style
# of row 180 (same even/odd parity as the new row 182).
$ws.Range("B180:C180").Copy()
$ws.Range("B" + $spacerRow + ":C" + $spacerRow).PasteSpecial(-4122)
$ws.Rows.Item($spacerRow).RowHeight = $ws.Rows.Item($lastRow).RowHeight

# --- Shift the existing data down by one row ------------------------------
# Capture the existing B/C text for rows 166-181 before overwriting
# anything.
$codes = @()
$names = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $codes += $ws.Cells.Item($r, 2).Text
    $names += $ws.Cells.Item($r, 3).Text
}

# Write them back shifted down by one row, starting from the bottom so a
# value is never overwritten before it has been captured.
for ($i = $codes.Length - 1; $i -ge 0; $i--) {
    $targetRow = $firstRow + 1 + $i
    $ws.Cells.Item($targetRow, 2).Value = $codes[$i]
    $ws.Cells.Item($targetRow, 3).Value = $names[$i]
}

# Row 166 now holds the newly inserted country.
$ws.Cells.Item($firstRow, 2).Value = "429"
$ws.Cells.Item($firstRow, 3).Value = "BELIZE"
